$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2" = 43
    "C3" = 27
    "C4" = 41
    "C5" = 55
    "C6" = 29
    "C7" = 53
    "C8" = 67
    "C9" = 79
    "C10" = 39
    "C11" = 24
    "C12" = 123
    "C13" = 51
    "C14" = 77
    "C15" = 32
    "C16" = 56
    "C17" = 40
    "C18" = 65
    "C19" = 63
    "C20" = 52
    "C21" = 57
    "C23" = 31
    "C25" = 58
    "C26" = 30
    "C27" = 97
    "C28" = 59
    "C29" = 74
    "C30" = 102
    "C32" = 73
    "C33" = 60
    "C34" = 19
    "C35" = 18
    "C36" = 98
    "C37" = 42
    "C38" = 72
    "C39" = 44
    "C40" = 85
    "C41" = 86
    "C42" = 70
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
